$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 2.269101333333333
$ws.Range("H2").Value = 6.807304
$ws.Range("I2").Value = 0.02891211995713196
$ws.Range("J2").Value = 0.02891211995713196
$ws.Range("M2").Value = 31.40242733333333
$ws.Range("N2").Value = 94.20728199999999
$ws.Range("O2").Value = 0.5334014788811394
$ws.Range("P2").Value = 0.5334014788811395
$ws.Range("Q2").Value = 71.25528973196978
$ws.Range("R2").Value = 641.297607587728
$ws.Range("S2").Value = 0.01542176754272309
$ws.Range("T2").Value = 0.0154217675427231
$ws.Range("G3").Value = 2.269101333333333
$ws.Range("H3").Value = 6.807304
$ws.Range("I3").Value = 0.02891211995713196
$ws.Range("J3").Value = 0.02891211995713196
$ws.Range("O3").Value = 0.327656036225058
$ws.Range("P3").Value = 0.327656036225058
$ws.Range("Q3").Value = 43.77045568493423
$ws.Range("R3").Value = 393.934101164408
$ws.Range("S3").Value = 0.009473230624017252
$ws.Range("T3").Value = 0.009473230624017252
$ws.Range("G4").Value = 2.269101333333333
$ws.Range("H4").Value = 6.807304
$ws.Range("I4").Value = 0.02891211995713196
$ws.Range("J4").Value = 0.02891211995713196
$ws.Range("M4").Value = 1.868202333333333
$ws.Range("N4").Value = 5.604607
$ws.Range("O4").Value = 0.03173327580290011
$ws.Range("P4").Value = 0.03173327580290011
$ws.Range("Q4").Value = 4.239140405503111
$ws.Range("R4").Value = 38.152263649528
$ws.Range("S4").Value = 0.0009174762766462008
$ws.Range("T4").Value = 0.0009174762766462008
$ws.Range("G5").Value = 2.269101333333333
$ws.Range("H5").Value = 6.807304
$ws.Range("I5").Value = 0.02891211995713196
$ws.Range("J5").Value = 0.02891211995713196
$ws.Range("M5").Value = 6.311623666666667
$ws.Range("N5").Value = 18.934871
$ws.Range("O5").Value = 0.1072092090909023
$ws.Range("P5").Value = 0.1072092090909024
$ws.Range("Q5").Value = 14.32171367753156
$ws.Range("R5").Value = 128.895423097784
$ws.Range("S5").Value = 0.00309964551374541
$ws.Range("T5").Value = 0.003099645513745411
$ws.Range("I6").Value = 0.7238963226334669
$ws.Range("J6").Value = 0.7238963226334669
$ws.Range("M6").Value = 31.40242733333333
$ws.Range("N6").Value = 94.20728199999999
$ws.Range("O6").Value = 0.5334014788811394
$ws.Range("P6").Value = 0.5334014788811395
$ws.Range("Q6").Value = 1784.076791381436
$ws.Range("R6").Value = 16056.69112243292
$ws.Range("S6").Value = 0.3861273690493097
$ws.Range("T6").Value = 0.3861273690493098
$ws.Range("I7").Value = 0.7238963226334669
$ws.Range("J7").Value = 0.7238963226334669
$ws.Range("O7").Value = 0.327656036225058
$ws.Range("P7").Value = 0.327656036225058
$ws.Range("S7").Value = 0.2371889997119775
$ws.Range("T7").Value = 0.2371889997119775
$ws.Range("I8").Value = 0.7238963226334669
$ws.Range("J8").Value = 0.7238963226334669
$ws.Range("M8").Value = 1.868202333333333
$ws.Range("N8").Value = 5.604607
$ws.Range("O8").Value = 0.03173327580290011
$ws.Range("P8").Value = 0.03173327580290011
$ws.Range("Q8").Value = 106.1388149751941
$ws.Range("R8").Value = 955.249334776747
$ws.Range("S8").Value = 0.02297160165883297
$ws.Range("T8").Value = 0.02297160165883297
$ws.Range("I9").Value = 0.7238963226334669
$ws.Range("J9").Value = 0.7238963226334669
$ws.Range("M9").Value = 6.311623666666667
$ws.Range("N9").Value = 18.934871
$ws.Range("O9").Value = 0.1072092090909023
$ws.Range("P9").Value = 0.1072092090909024
$ws.Range("Q9").Value = 358.5844234302546
$ws.Range("R9").Value = 3227.259810872291
$ws.Range("S9").Value = 0.07760835221334665
$ws.Range("T9").Value = 0.07760835221334667
$ws.Range("G10").Value = 18.57257166666666
$ws.Range("H10").Value = 55.717715
$ws.Range("I10").Value = 0.2366454120188096
$ws.Range("J10").Value = 0.2366454120188096
$ws.Range("M10").Value = 31.40242733333333
$ws.Range("N10").Value = 94.20728199999999
$ws.Range("O10").Value = 0.5334014788811394
$ws.Range("P10").Value = 0.5334014788811395
$ws.Range("Q10").Value = 583.2238321556255
$ws.Range("R10").Value = 5249.01448940063
$ws.Range("S10").Value = 0.1262270127412696
$ws.Range("T10").Value = 0.1262270127412696
$ws.Range("G11").Value = 18.57257166666666
$ws.Range("H11").Value = 55.717715
$ws.Range("I11").Value = 0.2366454120188096
$ws.Range("J11").Value = 0.2366454120188096
$ws.Range("O11").Value = 0.327656036225058
$ws.Range("P11").Value = 0.327656036225058
$ws.Range("Q11").Value = 358.2607410030894
$ws.Range("R11").Value = 3224.346669027805
$ws.Range("S11").Value = 0.07753829769292885
$ws.Range("T11").Value = 0.07753829769292886
$ws.Range("G12").Value = 18.57257166666666
$ws.Range("H12").Value = 55.717715
$ws.Range("I12").Value = 0.2366454120188096
$ws.Range("J12").Value = 0.2366454120188096
$ws.Range("M12").Value = 1.868202333333333
$ws.Range("N12").Value = 5.604607
$ws.Range("O12").Value = 0.03173327580290011
$ws.Range("P12").Value = 0.03173327580290011
$ws.Range("Q12").Value = 34.69732172366722
$ws.Range("R12").Value = 312.275895513005
$ws.Range("S12").Value = 0.007509534127083816
$ws.Range("T12").Value = 0.007509534127083817
$ws.Range("G13").Value = 18.57257166666666
$ws.Range("H13").Value = 55.717715
$ws.Range("I13").Value = 0.2366454120188096
$ws.Range("J13").Value = 0.2366454120188096
$ws.Range("M13").Value = 6.311623666666667
$ws.Range("N13").Value = 18.934871
$ws.Range("O13").Value = 0.1072092090909023
$ws.Range("P13").Value = 0.1072092090909024
$ws.Range("Q13").Value = 117.2230828821961
$ws.Range("R13").Value = 1055.007745939765
$ws.Range("S13").Value = 0.02537056745752729
$ws.Range("T13").Value = 0.02537056745752729
$ws.Range("G14").Value = 0.8276899999999999
$ws.Range("H14").Value = 2.48307
$ws.Range("I14").Value = 0.01054614539059158
$ws.Range("J14").Value = 0.01054614539059158
$ws.Range("M14").Value = 31.40242733333333
$ws.Range("N14").Value = 94.20728199999999
$ws.Range("O14").Value = 0.5334014788811394
$ws.Range("P14").Value = 0.5334014788811395
$ws.Range("Q14").Value = 25.99147507952666
$ws.Range("R14").Value = 233.9232757157399
$ws.Range("S14").Value = 0.005625329547837061
$ws.Range("T14").Value = 0.005625329547837062
$ws.Range("G15").Value = 0.8276899999999999
$ws.Range("H15").Value = 2.48307
$ws.Range("I15").Value = 0.01054614539059158
$ws.Range("J15").Value = 0.01054614539059158
$ws.Range("O15").Value = 0.327656036225058
$ws.Range("P15").Value = 0.327656036225058
$ws.Range("Q15").Value = 15.96595442154333
$ws.Range("R15").Value = 143.69358979389
$ws.Range("S15").Value = 0.003455508196134404
$ws.Range("T15").Value = 0.003455508196134404
$ws.Range("G16").Value = 0.8276899999999999
$ws.Range("H16").Value = 2.48307
$ws.Range("I16").Value = 0.01054614539059158
$ws.Range("J16").Value = 0.01054614539059158
$ws.Range("M16").Value = 1.868202333333333
$ws.Range("N16").Value = 5.604607
$ws.Range("O16").Value = 0.03173327580290011
$ws.Range("P16").Value = 0.03173327580290011
$ws.Range("Q16").Value = 1.546292389276666
$ws.Range("R16").Value = 13.91663150349
$ws.Range("S16").Value = 0.0003346637403371263
$ws.Range("T16").Value = 0.0003346637403371263
$ws.Range("G17").Value = 0.8276899999999999
$ws.Range("H17").Value = 2.48307
$ws.Range("I17").Value = 0.01054614539059158
$ws.Range("J17").Value = 0.01054614539059158
$ws.Range("M17").Value = 6.311623666666667
$ws.Range("N17").Value = 18.934871
$ws.Range("O17").Value = 0.1072092090909023
$ws.Range("P17").Value = 0.1072092090909024
$ws.Range("Q17").Value = 5.224067792663333
$ws.Range("R17").Value = 47.01661013397
$ws.Range("S17").Value = 0.001130643906282989
$ws.Range("T17").Value = 0.001130643906282989
